$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.366.84"
$ws.Range("E2").Value = "  +2.81%  "

$ws.Range("D3").Value = "1.579.79"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  +1.35%  "

$ws.Range("D5").Value = "'212.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "

$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("E7").Value = "  +1.27%  "

$ws.Range("D8").Value = "'46.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.16%  "

$ws.Range("D9").Value = "'24.02"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.54%  "

$ws.Range("E10").Value = "  -0.41%  "

$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("E12").Value = "  +1.16%  "

$ws.Range("D13").Value = "1.804.13"
$ws.Range("E13").Value = "  +0.52%  "

$ws.Range("D14").Value = "1.572.71"
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("E15").Value = "  +1.18%  "

$ws.Range("E16").Value = "  -0.95%  "

$ws.Range("D17").Value = "28.388.50"
$ws.Range("E17").Value = "  +3.00%  "

$ws.Range("D18").Value = "'62.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.58%  "

$ws.Range("D19").Value = "'229.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("E20").Value = "  -0.42%  "

$ws.Range("E21").Value = "  -1.12%  "

$ws.Range("E22").Value = "  +1.33%  "

$ws.Range("D23").Value = "'3.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.60%  "

$ws.Range("D24").Value = "'9.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.80%  "

$ws.Range("E25").Value = "  +3.82%  "

$ws.Range("D26").Value = "'151.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.01%  "

$ws.Range("D27").Value = "'15.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.09%  "

$ws.Range("D28").Value = "'6.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.12%  "

$ws.Range("E29").Value = "  -1.59%  "

$ws.Range("E30").Value = "  +1.23%  "

$ws.Range("E32").Value = "  -1.19%  "

$ws.Range("E33").Value = "  -0.28%  "

$ws.Range("D34").Value = "'3.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.89%  "

$ws.Range("D35").Value = "1.392.02"
$ws.Range("E35").Value = "  -4.48%  "

$ws.Range("D36").Value = "'1.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.37%  "

$ws.Range("D37").Value = "'1.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.17%  "

$ws.Range("D38").Value = "'2.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.35%  "

$ws.Range("E39").Value = "  +6.24%  "

$ws.Range("E40").Value = "  -1.11%  "

$ws.Range("D41").Value = "'0.536"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("D42").Value = "'0.806"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.72%  "

$ws.Range("E43").Value = "  +1.29%  "

$ws.Range("D44").Value = "'1.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.74%  "

$ws.Range("E45").Value = "  -1.26%  "

$ws.Range("D46").Value = "'0.984"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.35%  "

$ws.Range("D47").Value = "'62.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.94%  "

$ws.Range("D48").Value = "1.714.34"
$ws.Range("E48").Value = "  +0.47%  "

$ws.Range("D49").Value = "'86.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.90%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0102"
$ws.Range("E50").Value = "  +4.16%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0518"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.89%  "
